$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SearchTest")

# Update A2 with the new value, replacing old "Fortis" text
$ws.Range("A2").Value = "Gary Abbott"

# Remove the now-stale rows A3 and A4 (previously "ati" and "Ash")
$ws.Range("A3").Value = $null
$ws.Range("A4").Value = $null

# Update the active selection to A2
$ws.Range("A2").Select()
